$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) for the extra call-summary columns ---
# Copy the format of the existing header (A1) onto the new header range so
# the new headers pick up the same bold/centered/bordered style, then set
# their text.
$ws.Range("A1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)
$ws.Range("D1").Value = "phone"
$ws.Range("E1").Value = "product"
$ws.Range("F1").Value = "call_status"
$ws.Range("G1").Value = "timestamp"

# --- Existing phone columns become numeric values instead of text ---
$ws.Range("A2").Value = 919510038048
$ws.Range("B2").Value = 919510038048
$ws.Range("A3").Value = 919328027733
$ws.Range("B3").Value = 919328027733
$ws.Range("A4").Value = 919106284482
$ws.Range("B4").Value = 919106284482

# --- New row 5 holding the call-summary record ---
# D5 looks numeric ("+919510038048") so a leading apostrophe keeps it text,
# matching how Excel stores typed-in text that resembles a number; then we
# reset the cell format back to the default (General) style picked up from
# an existing plain cell, since the quote-prefix alone would otherwise tag
# the cell with its own style.
$ws.Range("D5").Value = "'+919510038048"
$ws.Range("C2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("E5").Value = "Laptop Pro "
$ws.Range("F5").Value = "Completed"
$ws.Range("G5").Value = "2025-11-27 16:14:47"
